$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1578, 1578, 1578, 0.01050180594126384),
    @(1, 911, 911, 911, 0.01052077611287435),
    @(2, 1233, 1233, 1233, 0.01263062953948975),
    @(3, 1162, 1162, 1162, 0.01228866577148438),
    @(4, 1356, 1356, 1356, 0.01285211245218913),
    @(5, 1484, 1484, 1484, 0.01248730818430583),
    @(6, 1921, 1921, 1921, 0.01237963835398356),
    @(7, 1748, 1748, 1748, 0.01079978148142497),
    @(8, 1645, 1645, 1645, 0.01245078245798747),
    @(9, 1741, 1741, 1741, 0.01251538594563802)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
